$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I0 and IF, styled like the existing header cells (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: I column, J column (J mirrors H column values)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 5

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 3

$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 8
